$d = $word.ActiveDocument

# --- 1. Remove the stray _GoBack bookmark from the first paragraph ---
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# --- 2. Merge the split run "...cativeiro" + "s ate..." into a single run ---
$d.Content.Find.Execute(
    "cativeiros até que alguém compre?", $false, $false, $false, $false, $false,
    $true, 1, $false, "cativeiros até que alguém compre?", 2) | Out-Null

# --- 3. Append the new "Android" section as a block of paragraphs, replacing the
#        final (empty) paragraph so the new content inherits its trailing position ---
$lastOld = $d.Paragraphs.Last
$insertRange = $d.Range($lastOld.Range.Start, $lastOld.Range.End)

$xmlFrag = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr></w:p>
<w:p><w:pPr><w:rPr><w:b/><w:color w:val="70AD47" w:themeColor="accent6"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:color w:val="70AD47" w:themeColor="accent6"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Android</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:color w:val="70AD47" w:themeColor="accent6"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>:</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">Vou falar sobre os métodos de desenvolvimento que vamos utilizar. </w:t></w:r><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">O aplicativo será desenvolvido para dispositivos </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Android</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>, sistema operacional móvel da Google. Hoje, é com folga o sistema móvel mais utilizado do mundo, sendo utilizado por aparelhos da Samsung, Motorola, HTC, entre outras gigantes no mercado de dispositivos móveis.</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:b/><w:color w:val="70AD47" w:themeColor="accent6"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:color w:val="70AD47" w:themeColor="accent6"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Android</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:color w:val="70AD47" w:themeColor="accent6"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> Studio:</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">Vamos usar como ambiente de desenvolvimento o IDE oficial da Google para desenvolvimento </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Android</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">, o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Android</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> Studio, que foi lançado em 2013</w:t></w:r><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">, desbancando o eclipse por já vir com o kit de desenvolvimento do </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Android</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> junto com o software.</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:b/><w:color w:val="70AD47" w:themeColor="accent6"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:color w:val="70AD47" w:themeColor="accent6"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Java:</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">E como toda aplicação </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Android</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>, o desenvolvimento será feito utilizando a linguagem de programação Java, uma das linguagens mais populares entre os desenvolvedores no mundo.</w:t></w:r></w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
"@

$insertRange.InsertXML($xmlFrag)

# InsertXML keeps the original final paragraph mark as a trailing empty paragraph;
# merge it away so the new "Java" paragraph becomes the true last paragraph.
$finalPara = $d.Paragraphs.Last
$newLastContentPara = $d.Paragraphs($d.Paragraphs.Count - 1)
$mergeRange = $d.Range($newLastContentPara.Range.End - 1, $finalPara.Range.End)
$mergeRange.Delete()

# --- 4. Re-add the _GoBack bookmark at the very end of the document ---
$endPoint = $d.Content.End - 1
$d.Bookmarks.Add("_GoBack", $d.Range($endPoint, $endPoint)) | Out-Null

Write-Output "OK"
Write-Output $d.Paragraphs.Count
